$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move completed cards from the "Text Cards to implement" column (C)
# into the "Finished Cards" column (D): Adventurer, Bureaucrat, Library.

# Adventurer: was C9 -> now D10
$ws.Range("C9").ClearContents() | Out-Null
$ws.Range("D10").Value = "Adventurer"

# Bureaucrat: was C10 -> now D11
$ws.Range("C10").ClearContents() | Out-Null
$ws.Range("D11").Value = "Bureaucrat"

# Library: was C17 -> now D21
$ws.Range("C17").ClearContents() | Out-Null
$ws.Range("D21").Value = "Library"

# Update the active selection to reflect where the user left off.
$ws.Range("C27").Select() | Out-Null
